# Perbaikan lebar tampilan profil + update data kendaraan dan user
# Update data_kendaraan sheet:
#   - row 3 (Rahma): Pajak_Terhutang / Tanggal_Jatuh_Tempo / Pajak become
#     plain text values instead of numbers/dates.
#   - new row 4 added for Siti Aminah (BG8989HI / Honda / Sepeda Motor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Add the new row 4 FIRST (while row 3's original date style is still
#    intact), so the date-format style slot used by F3 can be reused for
#    F4's numeric date value.
# ---------------------------------------------------------------------

$a4 = $ws.Cells.Item(4, 1)
$a4.NumberFormat = "@"
$a4.Value2 = "7868866666665555"
$a4.Style = "Normal"

$b4 = $ws.Cells.Item(4, 2)
$b4.Value2 = "BG8989HI"

$c4 = $ws.Cells.Item(4, 3)
$c4.Value2 = "Siti Aminah"

$d4 = $ws.Cells.Item(4, 4)
$d4.Value2 = "Palembang"

$e4 = $ws.Cells.Item(4, 5)
$e4.Value2 = 65000

$f4 = $ws.Cells.Item(4, 6)
$f4.NumberFormat = "YYYY-MM-DD"
$f4.Value2 = 46240

$g4 = $ws.Cells.Item(4, 7)
$g4.Value2 = 65000

$h4 = $ws.Cells.Item(4, 8)
$h4.Value2 = "MKK89JKK"

$i4 = $ws.Cells.Item(4, 9)
$i4.Value2 = "Honda"

$j4 = $ws.Cells.Item(4, 10)
$j4.Value2 = "Sepeda Motor"

$k4 = $ws.Cells.Item(4, 11)
$k4.Value2 = "Hitam"

foreach ($col in 12, 13, 14, 15) {
    $blank4 = $ws.Cells.Item(4, $col)
    $blank4.NumberFormat = "@"
    $blank4.Style = "Normal"
}

# ---------------------------------------------------------------------
# 2) Convert row 3's Pajak_Terhutang / Tanggal_Jatuh_Tempo / Pajak cells
#    from numbers/dates into plain text, matching the new text-based
#    format used for this record.
# ---------------------------------------------------------------------

$e3 = $ws.Cells.Item(3, 5)
$e3.NumberFormat = "@"
$e3.Value2 = "65000"
$e3.Style = "Normal"

$f3 = $ws.Cells.Item(3, 6)
$f3.NumberFormat = "@"
$f3.Value2 = "2026-08-06 00:00:00"
$f3.Style = "Normal"

$g3 = $ws.Cells.Item(3, 7)
$g3.NumberFormat = "@"
$g3.Value2 = "65000"
$g3.Style = "Normal"
